# SBA Science Variable Dictionary 2018-2019.xlsx
#
# The original workbook has a single sheet "SBAScience2019" that ends with
# a "Last edited on 7/11/2019" note in cell I2. This edit duplicates that
# sheet into a new sheet "SBASciFALL1819" (placed right after the original),
# makes the new sheet the active tab, and bumps the "last edited" note on
# the new copy to "Last edited on 8/14/2019" - the original sheet's note is
# left untouched.

$wb = $excel.ActiveWorkbook

# The existing (only) worksheet.
$ws1 = $wb.Worksheets.Item(1)

# Duplicate it immediately after itself - this becomes sheet 2 / "SBASciFALL1819".
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "SBASciFALL1819"

# Update the "last edited" note on the new copy only.
$ws2.Range("I2").Value = "Last edited on 8/14/2019"

# Deselect everything on the original sheet (mirrors it no longer being the
# active/focused tab) by selecting the whole sheet there.
$ws1.Cells.Select()

# Make the new sheet the active tab, with I3 selected (just below the note).
$ws2.Activate()
$ws2.Range("I3").Select()
